$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) First paragraph: "This is a Microsoft word document." gets two
#    trailing spaces appended, followed by a new red (C00000) run
#    "(This is a change – Version for branch alternate)" split across
#    three runs (as produced by incremental typing/formatting).
# ------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$bodyRange = $d.Range($p1.Range.Start, $p1.Range.End - 1)
$bodyRange.Text = "This is a Microsoft word document.  "

$dash = [char]0x2013

$pos = $bodyRange.End
$chunk1 = "(This is a change " + $dash + " Ve"
$ip = $d.Range($pos, $pos)
$ip.InsertAfter($chunk1)
$run1 = $d.Range($pos, $pos + $chunk1.Length)
$run1.Font.Color = 192

$pos = $run1.End
$chunk2 = "rsion for branch alternate"
$ip = $d.Range($pos, $pos)
$ip.InsertAfter($chunk2)
$run2 = $d.Range($pos, $pos + $chunk2.Length)
$run2.Font.Color = 192

$pos = $run2.End
$chunk3 = ")"
$ip = $d.Range($pos, $pos)
$ip.InsertAfter($chunk3)
$run3 = $d.Range($pos, $pos + $chunk3.Length)
$run3.Font.Color = 192

# ------------------------------------------------------------------
# 2) Insert a new, empty paragraph right after the "It will be
#    treated as a binary file by Git." paragraph, carrying shading
#    (fill F9F9F9) and bold Calibri/Times New Roman paragraph mark
#    formatting (color 202122), and nothing else.
# ------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(2)

$newParaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="F9F9F9"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:b/><w:bCs/><w:color w:val="202122"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insAfter = $d.Range($p2.Range.End, $p2.Range.End)
$insAfter.InsertXML($newParaXml)
